$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply scraped-value updates cell by cell. Column D (Price) values are
# prefixed with a leading apostrophe so Excel stores them as text exactly
# as scraped (preserving values like "1.00" or "61.891.56" rather than
# auto-converting them to numbers).

$ws.Range("D2").Value = "'61.891.56"
$ws.Range("E2").Value = '  -3.39%  '

$ws.Range("D3").Value = "'3.019.90"
$ws.Range("E3").Value = '  -4.07%  '

$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").Value = "'527.60"
$ws.Range("E5").Value = '  -5.92%  '

$ws.Range("D6").Value = "'128.34"
$ws.Range("E6").Value = '  -8.75%  '

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").Value = "'3.008.86"
$ws.Range("E8").Value = '  -4.23%  '

$ws.Range("E9").Value = '  -0.20%  '

$ws.Range("D10").Value = "'0.149"
$ws.Range("E10").Value = '  -2.63%  '

$ws.Range("D11").Value = "'6.02"
$ws.Range("E11").Value = '  -10.14%  '

$ws.Range("D12").Value = "'0.441"
$ws.Range("E12").Value = '  -4.51%  '

$ws.Range("D13").Value = "'0.0000220"
$ws.Range("E13").Value = '  +0.34%  '

$ws.Range("D14").Value = "'33.03"
$ws.Range("E14").Value = '  -8.79%  '

$ws.Range("D15").Value = "'3.508.99"
$ws.Range("E15").Value = '  -4.13%  '

$ws.Range("D16").Value = "'61.858.75"
$ws.Range("E16").Value = '  -3.62%  '

$ws.Range("E17").Value = '  -2.54%  '

$ws.Range("D18").Value = "'3.018.88"

$ws.Range("E19").Value = '  -4.97%  '

$ws.Range("D20").Value = "'470.38"
$ws.Range("E20").Value = '  -7.67%  '

$ws.Range("D21").Value = "'12.93"
$ws.Range("E21").Value = '  -7.17%  '

$ws.Range("E22").Value = '  -4.55%  '

$ws.Range("D23").Value = "'6.81"
$ws.Range("E23").Value = '  -7.86%  '

$ws.Range("D24").Value = "'77.71"
$ws.Range("E24").Value = '  -0.96%  '

$ws.Range("D25").Value = "'11.67"
$ws.Range("E25").Value = '  -8.22%  '

$ws.Range("E26").Value = '  +0.09%  '

$ws.Range("E27").Value = '  -7.22%  '

$ws.Range("D28").Value = "'7.88"
$ws.Range("E28").Value = '  -9.27%  '

$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = '  -0.09%  '

$ws.Range("D30").Value = "'25.25"
$ws.Range("E30").Value = '  -4.70%  '

$ws.Range("D31").Value = "'1.80"
$ws.Range("E31").Value = '  -14.11%  '

$ws.Range("E32").Value = '  -4.81%  '

$ws.Range("D33").Value = "'56.29"
$ws.Range("E33").Value = '  +4.71%  '

$ws.Range("E34").Value = '  -10.84%  '

$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").Value = "'5.12"
$ws.Range("E35").Value = '  -3.47%  '

$ws.Range("B36").Value = 'Filecoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D36").Value = "'5.77"
$ws.Range("E36").Value = '  -4.31%  '

$ws.Range("D37").Value = "'462.00"
$ws.Range("E37").Value = '  -16.16%  '

$ws.Range("D38").Value = "'3.032.96"
$ws.Range("E38").Value = '  -3.79%  '

$ws.Range("D39").Value = "'0.0383"
$ws.Range("E39").Value = '  -9.85%  '

$ws.Range("E40").Value = '  -5.32%  '

$ws.Range("E41").Value = '  -8.07%  '

$ws.Range("D42").Value = "'7.88"
$ws.Range("E42").Value = '  -4.09%  '

$ws.Range("D43").Value = "'2.48"
$ws.Range("E43").Value = '  -8.95%  '

$ws.Range("E44").Value = '  +0.07%  '

$ws.Range("E45").Value = '  -7.66%  '

$ws.Range("D46").Value = "'0.0₃0521"
$ws.Range("E46").Value = '  +1.89%  '

$ws.Range("D47").Value = "'117.70"
$ws.Range("E47").Value = '  -3.42%  '

$ws.Range("E48").Value = '  -9.91%  '

$ws.Range("E49").Value = '  -1.55%  '

$ws.Range("D50").Value = "'23.61"
$ws.Range("E50").Value = '  -4.70%  '

$ws.Range("D51").Value = "'2.27"
$ws.Range("E51").Value = '  +2.91%  '

